# Update "想去人数" (want-to-go count) figures in column F across sheets
# 展览 (Exhibitions), 演出 (Performances) and 全部类型 (All types).
# Mirrors the regenerated stats published at gh-pages output 456a3b4.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 12880
$ws.Range("F5").Value = 38
$ws.Range("F6").Value = 326
$ws.Range("F7").Value = 407
$ws.Range("F9").Value = 12957
$ws.Range("F10").Value = 43
$ws.Range("F11").Value = 28
$ws.Range("F12").Value = 5272
$ws.Range("F13").Value = 549
$ws.Range("F14").Value = 21
$ws.Range("F15").Value = 15
$ws.Range("F16").Value = 33
$ws.Range("F18").Value = 41
$ws.Range("F19").Value = 135
$ws.Range("F20").Value = 684
$ws.Range("F22").Value = 6199
$ws.Range("F23").Value = 1163
$ws.Range("F24").Value = 3635
$ws.Range("F26").Value = 47

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 9

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 12880
$ws.Range("F5").Value = 38
$ws.Range("F6").Value = 326
$ws.Range("F8").Value = 407
$ws.Range("F10").Value = 12957
$ws.Range("F11").Value = 43
$ws.Range("F12").Value = 28
$ws.Range("F13").Value = 5272
$ws.Range("F14").Value = 549
$ws.Range("F15").Value = 21
$ws.Range("F16").Value = 15
$ws.Range("F17").Value = 33
$ws.Range("F19").Value = 41
$ws.Range("F20").Value = 135
$ws.Range("F21").Value = 684
$ws.Range("F23").Value = 9
$ws.Range("F24").Value = 6199
$ws.Range("F25").Value = 1163
$ws.Range("F26").Value = 3635
$ws.Range("F28").Value = 47
